# Updating simple excel file template
#
# Reproduces the content-level changes from the commit:
#   - registerinfo!B5 hyperlink text: .../test-bba/testreg1 -> .../utils/testreg1
#   - selection on registerinfo moves from B5 to B6 (and loses tab focus)
#   - testreg becomes the active/selected sheet (tabSelected / activeTab)
#   - testreg column widths are narrowed slightly (cols D:E end up equal width)

$wb  = $excel.ActiveWorkbook
$registerinfo = $wb.Worksheets.Item(1)   # sheetId=2, r:id=rId1, xl/worksheets/sheet1.xml
$testreg      = $wb.Worksheets.Item(2)   # sheetId=1, r:id=rId2, xl/worksheets/sheet2.xml

# --- sharedStrings.xml: update the registry_location hyperlink cell text ---
$registerinfo.Range("B5").Value = "http://registry.it.csiro.au/sandbox/csiro/utils/testreg1"

# --- testreg!cols: new column widths (D and E end up the same width) ---
$testreg.Columns.Item(1).ColumnWidth = 20.666666666666664
$testreg.Columns.Item(2).ColumnWidth = 14.666666666666668
$testreg.Columns.Item(3).ColumnWidth = 18.666666666666664
$testreg.Range($testreg.Columns.Item(4), $testreg.Columns.Item(5)).ColumnWidth = 17
$testreg.Columns.Item(6).ColumnWidth = 24

# --- sheet selections ---
# registerinfo: active cell moves from B5 to B6
$registerinfo.Range("B6").Select() | Out-Null

# testreg becomes the active sheet/tab (moves tabSelected + workbook activeTab)
$testreg.Activate() | Out-Null
